$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so values such as
# "77.90" or "0.999" are stored as literal strings (matching the
# original inline-string cells) instead of being auto-converted
# to numbers by Excel's usual value-entry parsing.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.833.82"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "2.356.84"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "239.86"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -1.99%  "

$ws.Range("D7").Value = "74.31"
$ws.Range("E7").Value = "  +0.85%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  +2.39%  "

$ws.Range("E10").Value = "  +1.77%  "

$ws.Range("D11").Value = "60.92"
$ws.Range("E11").Value = "  +6.33%  "

$ws.Range("D12").Value = "35.75"
$ws.Range("E12").Value = "  +11.49%  "

$ws.Range("E13").Value = "  +0.80%  "

$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("D17").Value = "2.359.83"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").Value = "43.793.48"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").Value = "77.90"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("E21").Value = "  -2.51%  "

$ws.Range("D22").Value = "251.10"
$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("D23").Value = "0.999"

$ws.Range("E24").Value = "  +3.32%  "

$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").Value = "2.50"
$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("D27").Value = "10.49"
$ws.Range("E27").Value = "  -2.33%  "

$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "22.33"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "175.52"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Value = "0.130"
$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("D32").Value = "0.133"
$ws.Range("E32").Value = "  -1.50%  "

$ws.Range("D33").Value = "0.0748"
$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("D34").Value = "5.07"
$ws.Range("E34").Value = "  -2.95%  "

$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").Value = "6.60"
$ws.Range("E37").Value = "  +5.01%  "

$ws.Range("D38").Value = "2.42"
$ws.Range("E38").Value = "  +2.90%  "

$ws.Range("D39").Value = "0.0279"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("E40").Value = "  +14.81%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "20.25"
$ws.Range("E41").Value = "  +7.03%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "65.30"
$ws.Range("E42").Value = "  +11.87%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.203"
$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("D44").Value = "9.04"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.106"
$ws.Range("E45").Value = "  -3.71%  "

$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("D49").Value = "1.16"
$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("D50").Value = "98.16"
$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("E51").Value = "  +2.17%  "

# Restore column D formatting back to the default "Normal" style so
# the cells end up with no explicit style index again, matching the
# original (unstyled) cells.
$dRange.Style = "Normal"

Write-Host "applied cryptos update"